# Refresh the crypto price/volume columns (D = Price, E = Volume(1h))
# to match the latest scrape, per the GitHub Actions update job.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    # Some "Price" values are plain numeric strings (e.g. "309.04").
    # Force the cell to Text first so Excel stores/keeps them as a
    # literal string (matching the source data) instead of silently
    # converting them to a number. Reset the style afterwards so we
    # do not leave a stray number-format on the cell.
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

$ws.Range('D2').Value = '27.844.03'
$ws.Range('E2').Value = '  -1.21%  '
$ws.Range('D3').Value = '1.810.32'
$ws.Range('E3').Value = '  +0.73%  '
$ws.Range('E4').Value = '  +0.01%  '
Set-TextValue 'D5' '309.04'
$ws.Range('E5').Value = '  -1.71%  '
$ws.Range('E6').Value = '  -0.03%  '
Set-TextValue 'D7' '0.4992'
$ws.Range('E7').Value = '  -4.14%  '
Set-TextValue 'D8' '0.3887'
$ws.Range('E8').Value = '  +1.75%  '
Set-TextValue 'D9' '0.09576'
$ws.Range('E9').Value = '  +20.66%  '
Set-TextValue 'D10' '1.098'
$ws.Range('E10').Value = '  -0.10%  '
Set-TextValue 'D11' '40.32'
$ws.Range('E11').Value = '  -2.79%  '
Set-TextValue 'D12' '6.404'
$ws.Range('E12').Value = '  +1.89%  '
Set-TextValue 'D13' '1.001'
$ws.Range('E13').Value = '  +0.01%  '
Set-TextValue 'D14' '20.45'
$ws.Range('E14').Value = '  -0.59%  '
$ws.Range('D15').Value = '1.816.62'
$ws.Range('E15').Value = '  +1.38%  '
Set-TextValue 'D16' '7.255'
$ws.Range('E16').Value = '  -0.04%  '
Set-TextValue 'D18' '93.37'
$ws.Range('E18').Value = '  +0.11%  '
Set-TextValue 'D19' '0.06595'
$ws.Range('E19').Value = '  +0.59%  '
Set-TextValue 'D20' '1.001'
$ws.Range('E20').Value = '  -0.01%  '
Set-TextValue 'D21' '17.14'
$ws.Range('E21').Value = '  -0.92%  '
Set-TextValue 'D22' '5.944'
$ws.Range('E22').Value = '  -0.17%  '
$ws.Range('D23').Value = '27.914.83'
$ws.Range('E23').Value = '  -1.11%  '
Set-TextValue 'D24' '11.17'
$ws.Range('E24').Value = '  +0.17%  '
Set-TextValue 'D25' '2.248'
$ws.Range('E25').Value = '  -0.88%  '
Set-TextValue 'D26' '157.67'
$ws.Range('E26').Value = '  -2.07%  '
Set-TextValue 'D27' '20.75'
$ws.Range('E27').Value = '  +1.43%  '
$ws.Range('D28').Value = '2.023.09'
$ws.Range('E28').Value = '  +1.13%  '
Set-TextValue 'D29' '2.400'
$ws.Range('E29').Value = '  +2.73%  '
Set-TextValue 'D30' '128.22'
$ws.Range('E30').Value = '  +4.08%  '
Set-TextValue 'D31' '0.1071'
$ws.Range('E31').Value = '  +0.36%  '
Set-TextValue 'D32' '1.053'
$ws.Range('E32').Value = '  -0.11%  '
Set-TextValue 'D33' '5.586'
$ws.Range('E33').Value = '  +0.34%  '
Set-TextValue 'D34' '3.623'
$ws.Range('E34').Value = '  -1.29%  '
Set-TextValue 'D35' '0.06807'
$ws.Range('E35').Value = '  -7.01%  '
Set-TextValue 'D36' '9.007'
$ws.Range('E36').Value = '  +4.47%  '
Set-TextValue 'D37' '0.02317'
$ws.Range('E37').Value = '  -0.43%  '
Set-TextValue 'D38' '0.2145'
$ws.Range('E38').Value = '  +0.14%  '
$ws.Range('E39').Value = '  -7.46%  '
Set-TextValue 'D40' '4.934'
$ws.Range('E40').Value = '  -2.74%  '
Set-TextValue 'D41' '0.6242'
$ws.Range('E41').Value = '  +1.28%  '
Set-TextValue 'D42' '1.001'
$ws.Range('E42').Value = '  +0.00%  '
$ws.Range('E43').Value = '  -1.59%  '
Set-TextValue 'D44' '13.07'
$ws.Range('E44').Value = '  -1.07%  '
Set-TextValue 'D45' '0.5916'
$ws.Range('E45').Value = '  -1.35%  '
Set-TextValue 'D46' '1.295'
$ws.Range('E46').Value = '  -5.36%  '
Set-TextValue 'D47' '3.687'
$ws.Range('E47').Value = '  -2.48%  '
Set-TextValue 'D48' '123.82'
$ws.Range('E48').Value = '  -2.91%  '
Set-TextValue 'D49' '1.957'
$ws.Range('E49').Value = '  +1.90%  '
Set-TextValue 'D50' '1.177'
$ws.Range('E50').Value = '  -4.35%  '
Set-TextValue 'D51' '0.06786'
$ws.Range('E51').Value = '  +0.11%  '
